$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '74.209.29'
$ws.Range('E2').Value = '  +7.71%  '
$ws.Range('D3').Value = '2.638.79'
$ws.Range('E3').Value = '  +7.94%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '186.25'
$ws.Range('E5').Value = '  +13.84%  '
$ws.Range('D6').Value = '582.75'
$ws.Range('E6').Value = '  +4.01%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.532'
$ws.Range('E8').Value = '  +4.73%  '
$ws.Range('D9').Value = '0.199'
$ws.Range('E9').Value = '  +17.51%  '
$ws.Range('D10').Value = '2.634.72'
$ws.Range('E10').Value = '  +7.83%  '
$ws.Range('E11').Value = '  +1.32%  '
$ws.Range('E12').Value = '  +8.05%  '
$ws.Range('D13').Value = '4.70'
$ws.Range('E13').Value = '  +1.89%  '
$ws.Range('E14').Value = '  +5.76%  '
$ws.Range('D15').Value = '74.225.34'
$ws.Range('E15').Value = '  +7.93%  '
$ws.Range('D16').Value = '3.120.44'
$ws.Range('E16').Value = '  +7.88%  '
$ws.Range('D17').Value = '26.39'
$ws.Range('E17').Value = '  +12.71%  '
$ws.Range('D18').Value = '2.644.46'
$ws.Range('E18').Value = '  +8.32%  '
$ws.Range('D19').Value = '9.12'
$ws.Range('E19').Value = '  +29.02%  '
$ws.Range('D20').Value = '11.85'
$ws.Range('E20').Value = '  +11.50%  '
$ws.Range('D21').Value = '371.90'
$ws.Range('E21').Value = '  +9.53%  '
$ws.Range('D22').Value = '2.29'
$ws.Range('E22').Value = '  +17.28%  '
$ws.Range('D23').Value = '4.07'
$ws.Range('E23').Value = '  +5.77%  '
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').Value = '69.73'
$ws.Range('E25').Value = '  +5.81%  '
$ws.Range('D26').Value = '4.13'
$ws.Range('E26').Value = '  +9.67%  '
$ws.Range('D27').Value = '9.30'
$ws.Range('E27').Value = '  +11.71%  '
$ws.Range('D28').Value = '2.778.03'
$ws.Range('E28').Value = '  +8.08%  '
$ws.Range('E29').Value = '  -0.63%  '
$ws.Range('D30').Value = '0.0₃0946'
$ws.Range('E30').Value = '  +14.53%  '
$ws.Range('D31').Value = '523.68'
$ws.Range('E31').Value = '  +20.94%  '
$ws.Range('E32').Value = '  +15.83%  '
$ws.Range('D33').Value = '7.67'
$ws.Range('E33').Value = '  +6.70%  '
$ws.Range('E34').Value = '  +9.35%  '
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').Value = '162.35'
$ws.Range('E36').Value = '  +1.35%  '
$ws.Range('D37').Value = '0.117'
$ws.Range('E37').Value = '  +10.15%  '
$ws.Range('D38').Value = '19.18'
$ws.Range('E38').Value = '  +6.30%  '
$ws.Range('E39').Value = '  +1.52%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').Value = '4.91'
$ws.Range('E41').Value = '  +11.80%  '
$ws.Range('D42').Value = '1.67'
$ws.Range('E42').Value = '  +9.88%  '
$ws.Range('D43').Value = '0.327'
$ws.Range('E43').Value = '  +8.63%  '
$ws.Range('D44').Value = '161.65'
$ws.Range('E44').Value = '  +23.26%  '
$ws.Range('D45').Value = '2.38'
$ws.Range('E45').Value = '  +14.40%  '
$ws.Range('D46').Value = '1.18'
$ws.Range('E46').Value = '  +9.23%  '
$ws.Range('D47').Value = '38.91'
$ws.Range('E47').Value = '  +3.75%  '
$ws.Range('E48').Value = '  +18.42%  '
$ws.Range('D49').Value = '3.61'
$ws.Range('E49').Value = '  +8.18%  '
$ws.Range('D50').Value = '0.530'
$ws.Range('E50').Value = '  +9.44%  '
$ws.Range('D51').Value = '20.79'
$ws.Range('E51').Value = '  +22.50%  '
